$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheets (1-based via Worksheets.Item):
#   1 = Personas
#   2 = Lugar
#   3 = Organizaciones
#   4 = Momento
# ----------------------------------------------------------------------
$wsPersonas = $wb.Worksheets.Item(1)
$wsLugar = $wb.Worksheets.Item(2)
$wsOrganizaciones = $wb.Worksheets.Item(3)
$wsMomento = $wb.Worksheets.Item(4)

# ----------------------------------------------------------------------
# Personas: fill in Jesus' brothers (rows 99-102, per98-per101) and note
# that Felipe (row 103, per102) is a brother of Herodes.
# ----------------------------------------------------------------------
$wsPersonas.Range("B99").Value = "Jacobo"
$wsPersonas.Range("D99").Value = "hermano de Jesús"

$wsPersonas.Range("B100").Value = "José"
$wsPersonas.Range("D100").Value = "hermano de Jesús"

$wsPersonas.Range("B101").Value = "Simón"
$wsPersonas.Range("D101").Value = "hermano de Jesús"

$wsPersonas.Range("B102").Value = "Judas"
$wsPersonas.Range("D102").Value = "hermano de Jesús"

$wsPersonas.Range("B103").Value = "Felipe"
$wsPersonas.Range("D103").Value = "hermano de Herodes"

# ----------------------------------------------------------------------
# Lugar: "tierra" (pla41, row42) gets a variación "mundo".
# ----------------------------------------------------------------------
$wsLugar.Range("C42").Value = "mundo"

# ----------------------------------------------------------------------
# Organizaciones:
#  - org3 "discípulos" (row4) gets a comentario "discípulos de Jesús"
#  - org6 "demonio"/"ángeles" (row7) comment moves from column E to D
#  - org19 (row20) gets a nombre "no creyentes"
# ----------------------------------------------------------------------
$wsOrganizaciones.Range("D4").Value = "discípulos de Jesús"

$wsOrganizaciones.Range("D7").Value = $wsOrganizaciones.Range("E7").Text
$wsOrganizaciones.Range("E7").Value = $null

$wsOrganizaciones.Range("B20").Value = "no creyentes"

# ----------------------------------------------------------------------
# Selections / active sheet: Personas becomes the active tab, with the
# other sheets sharing the same "B99:D102" multi-area selection (plus a
# sheet specific cell) that the edit left behind when fixing the rows
# above.
# ----------------------------------------------------------------------
$wsPersonas.Activate()
$wsPersonas.Range("B99:D102").Select()

$wsLugar.Range("B58").Select()

$wsOrganizaciones.Range("B9").Select()

$wsMomento.Range("B3").Select()

$wsPersonas.Activate()
